# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Brynhildr_Profits workbook (profit calc columns H..N per leve row).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 52.75
$ws.Range("I6").Value = 52.75
$ws.Range("K6").Value = 158.25
$ws.Range("M6").Value = -46.25

$ws.Range("H12").Value = 159.45454
$ws.Range("I12").Value = 159.45454
$ws.Range("K12").Value = 159.45454
$ws.Range("M12").Value = 10.54545999999999

$ws.Range("H15").Value = 983.1111
$ws.Range("I15").Value = 983.1111
$ws.Range("K15").Value = 2949.3333
$ws.Range("M15").Value = -2780.3333

$ws.Range("H21").Value = 4573
$ws.Range("I21").Value = 5850
$ws.Range("K21").Value = 5850
$ws.Range("M21").Value = -5382

$ws.Range("H23").Value = 4573
$ws.Range("I23").Value = 5850
$ws.Range("K23").Value = 5850
$ws.Range("M23").Value = -5616

$ws.Range("H29").Value = 4194.4443
$ws.Range("J29").Value = 4992.857
$ws.Range("L29").Value = 14978.571
$ws.Range("N29").Value = -15540.571

$ws.Range("H38").Value = 9
$ws.Range("I38").Value = 9
$ws.Range("K38").Value = 27
$ws.Range("M38").Value = 345

$ws.Range("H58").Value = 2532.889
$ws.Range("J58").Value = 2749.3333
$ws.Range("L58").Value = 8247.999899999999
$ws.Range("N58").Value = -8547.999899999999

$ws.Range("H69").Value = 9946.5
$ws.Range("I69").Value = 8936
$ws.Range("J69").Value = 14999
$ws.Range("K69").Value = 26808
$ws.Range("L69").Value = 44997
$ws.Range("M69").Value = -25934
$ws.Range("N69").Value = -46745

$ws.Range("H72").Value = 9946.5
$ws.Range("I72").Value = 8936
$ws.Range("J72").Value = 14999
$ws.Range("K72").Value = 80424
$ws.Range("L72").Value = 134991
$ws.Range("M72").Value = -76056
$ws.Range("N72").Value = -143727

$ws.Range("H87").Value = 67800
$ws.Range("J87").Value = 67800
$ws.Range("L87").Value = 67800
$ws.Range("N87").Value = -70296

$ws.Range("H90").Value = 67800
$ws.Range("J90").Value = 67800
$ws.Range("L90").Value = 203400
$ws.Range("N90").Value = -215880

$ws.Range("H135").Value = 1800.9615
$ws.Range("I135").Value = 1666.1818
$ws.Range("J135").Value = 2542.25
$ws.Range("K135").Value = 14995.6362
$ws.Range("L135").Value = 22880.25
$ws.Range("M135").Value = -12460.6362
$ws.Range("N135").Value = -27950.25

$ws.Range("H137").Value = 3636.3408
$ws.Range("I137").Value = 1381.1471
$ws.Range("J137").Value = 11304
$ws.Range("K137").Value = 4143.4413
$ws.Range("L137").Value = 33912
$ws.Range("M137").Value = -1593.4413
$ws.Range("N137").Value = -39012

$ws.Range("H138").Value = 3961.1636
$ws.Range("J138").Value = 4506.3335
$ws.Range("L138").Value = 13519.0005
$ws.Range("N138").Value = -23799.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2999
$ws.Range("I61").Value = 2999
$ws.Range("K61").Value = 2999
$ws.Range("M61").Value = -2787

$ws.Range("H122").Value = 1609
$ws.Range("I122").Value = 1609
$ws.Range("K122").Value = 4827
$ws.Range("M122").Value = -2377

$ws.Range("H132").Value = 4055.8215
$ws.Range("I132").Value = 2941.75
$ws.Range("J132").Value = 5541.25
$ws.Range("K132").Value = 8825.25
$ws.Range("L132").Value = 16623.75
$ws.Range("M132").Value = -6295.25
$ws.Range("N132").Value = -21683.75

$ws.Range("H133").Value = 50750.285
$ws.Range("J133").Value = 50750.285
$ws.Range("L133").Value = 50750.285
$ws.Range("N133").Value = -55810.285

$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 2999
$ws.Range("K136").Value = 8997
$ws.Range("M136").Value = -6447

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2370.4285
$ws.Range("I86").Value = 2175.625
$ws.Range("J86").Value = 2993.8
$ws.Range("K86").Value = 2175.625
$ws.Range("L86").Value = 2993.8
$ws.Range("M86").Value = -1052.625
$ws.Range("N86").Value = -5239.8

$ws.Range("H89").Value = 2370.4285
$ws.Range("I89").Value = 2175.625
$ws.Range("J89").Value = 2993.8
$ws.Range("K89").Value = 10878.125
$ws.Range("L89").Value = 14969
$ws.Range("M89").Value = -5262.125
$ws.Range("N89").Value = -26201

$ws.Range("H111").Value = 39666.668
$ws.Range("J111").Value = 39666.668
$ws.Range("L111").Value = 39666.668
$ws.Range("N111").Value = -47846.668

$ws.Range("H132").Value = 79965
$ws.Range("I132").Value = 76709
$ws.Range("K132").Value = 76709
$ws.Range("M132").Value = -71649

$ws.Range("H134").Value = 4412.579
$ws.Range("I134").Value = 2299.2744
$ws.Range("J134").Value = 22375.666
$ws.Range("K134").Value = 6897.823199999999
$ws.Range("L134").Value = 67126.99800000001
$ws.Range("M134").Value = -4362.823199999999
$ws.Range("N134").Value = -72196.99800000001

$ws.Range("H141").Value = 243332.67
$ws.Range("J141").Value = 243332.67
$ws.Range("L141").Value = 243332.67
$ws.Range("N141").Value = -253692.67

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7245.3096
$ws.Range("I58").Value = 4380.189
$ws.Range("J58").Value = 28447.2
$ws.Range("K58").Value = 4380.189
$ws.Range("L58").Value = 28447.2
$ws.Range("M58").Value = -4177.189
$ws.Range("N58").Value = -28853.2

$ws.Range("H86").Value = 39497.742
$ws.Range("I86").Value = 57614.277
$ws.Range("K86").Value = 57614.277
$ws.Range("M86").Value = -56491.277

$ws.Range("H89").Value = 39497.742
$ws.Range("I89").Value = 57614.277
$ws.Range("K89").Value = 288071.385
$ws.Range("M89").Value = -282455.385

$ws.Range("H99").Value = 21032.545
$ws.Range("I99").Value = 35643.168
$ws.Range("J99").Value = 3499.8
$ws.Range("K99").Value = 35643.168
$ws.Range("L99").Value = 3499.8
$ws.Range("M99").Value = -34145.168
$ws.Range("N99").Value = -6495.8

$ws.Range("H126").Value = 21032.545
$ws.Range("I126").Value = 35643.168
$ws.Range("J126").Value = 3499.8
$ws.Range("K126").Value = 106929.504
$ws.Range("L126").Value = 10499.4
$ws.Range("M126").Value = -104459.504
$ws.Range("N126").Value = -15439.4

$ws.Range("H134").Value = 1189.3798
$ws.Range("I134").Value = 1198.7307
$ws.Range("K134").Value = 3596.1921
$ws.Range("M134").Value = -1061.1921

$ws.Range("H136").Value = 7245.3096
$ws.Range("I136").Value = 4380.189
$ws.Range("J136").Value = 28447.2
$ws.Range("K136").Value = 13140.567
$ws.Range("L136").Value = 85341.60000000001
$ws.Range("M136").Value = -10590.567
$ws.Range("N136").Value = -90441.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 909.0769
$ws.Range("I50").Value = 523.55554
$ws.Range("K50").Value = 1570.66662
$ws.Range("M50").Value = -1089.66662

$ws.Range("H53").Value = 909.0769
$ws.Range("I53").Value = 523.55554
$ws.Range("K53").Value = 1570.66662
$ws.Range("M53").Value = -1089.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1874.6666
$ws.Range("I43").Value = 1874.6666
$ws.Range("K43").Value = 1874.6666
$ws.Range("M43").Value = -1723.6666

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H69").Value = 33000
$ws.Range("J69").Value = 32000
$ws.Range("L69").Value = 32000
$ws.Range("N69").Value = -33498

$ws.Range("H72").Value = 33000
$ws.Range("J72").Value = 32000
$ws.Range("L72").Value = 96000
$ws.Range("N72").Value = -103488

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H126").Value = 4133.3335
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 4133.3335
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 12400.0005
$ws.Range("N126").Value = -17340.0005
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 18939.24
$ws.Range("I132").Value = 21755.572
$ws.Range("K132").Value = 65266.716
$ws.Range("M132").Value = -62736.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3254.818
$ws.Range("I22").Value = 1536.1666
$ws.Range("K22").Value = 1536.1666
$ws.Range("M22").Value = -1241.1666

$ws.Range("H27").Value = 3254.818
$ws.Range("I27").Value = 1536.1666
$ws.Range("K27").Value = 1536.1666
$ws.Range("M27").Value = -1429.1666

$ws.Range("H82").Value = 5685.174
$ws.Range("I82").Value = 4629.4287
$ws.Range("K82").Value = 4629.4287
$ws.Range("M82").Value = -4268.4287

$ws.Range("H85").Value = 5685.174
$ws.Range("I85").Value = 4629.4287
$ws.Range("K85").Value = 4629.4287
$ws.Range("M85").Value = -3381.4287

$ws.Range("H108").Value = 52900
$ws.Range("J108").Value = 52900
$ws.Range("L108").Value = 52900
$ws.Range("N108").Value = -60580

$ws.Range("H132").Value = 1825
$ws.Range("I132").Value = 1825
$ws.Range("K132").Value = 5475
$ws.Range("M132").Value = -2945

$ws.Range("H136").Value = 3162.24
$ws.Range("I136").Value = 3104.389
$ws.Range("K136").Value = 9313.167000000001
$ws.Range("M136").Value = -6763.167000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1324.3948
$ws.Range("I136").Value = 1286.0385
$ws.Range("K136").Value = 3858.1155
$ws.Range("M136").Value = -1308.1155
